$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.994.46'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '2.755.33'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.40'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.43'
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -3.54%  '
$ws.Range('E9').Value = '  -1.95%  '
$ws.Range('E10').Value = '  +2.40%  '
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.64'
$ws.Range('E12').Value = '  -16.56%  '
$ws.Range('D13').Value = '3.242.82'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.88'
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('D15').Value = '63.592.64'
$ws.Range('E15').Value = '  -0.69%  '
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = '2.758.11'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.18'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '356.55'
$ws.Range('E20').Value = '  -2.01%  '
$ws.Range('E21').Value = '  -3.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  -3.25%  '
$ws.Range('E24').Value = '  -2.58%  '
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.62'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').Value = '  -1.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.96'
$ws.Range('E29').Value = '  -2.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.18'
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.27'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '170.19'
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.95'
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.23'
$ws.Range('E34').Value = '  -2.07%  '
$ws.Range('E35').Value = '  +0.94%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('E39').Value = '  +1.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '336.50'
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.20'
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.16'
$ws.Range('E42').Value = '  -0.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.46'
$ws.Range('E43').Value = '  -1.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.74'
$ws.Range('E44').Value = '  -2.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0588'
$ws.Range('E45').Value = '  -2.40%  '
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('E48').Value = '  -2.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '135.11'
$ws.Range('E49').Value = '  -1.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.05'
$ws.Range('E51').Value = '  +0.02%  '
